$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.141.49'
$ws.Range("E2").Value = '  +0.70%  '

$ws.Range("D3").Value = '1.895.95'
$ws.Range("E3").Value = '  +0.63%  '

$ws.Range("D4").Value = '''1.004'
$ws.Range("E4").Value = '  +0.32%  '

$ws.Range("D5").Value = '''323.24'
$ws.Range("E5").Value = '  -1.96%  '

$ws.Range("D6").Value = '''1.004'
$ws.Range("E6").Value = '  +0.37%  '

$ws.Range("D7").Value = '''0.4702'
$ws.Range("E7").Value = '  +2.58%  '

$ws.Range("E8").Value = '  -2.56%  '

$ws.Range("D9").Value = '''47.44'
$ws.Range("E9").Value = '  -0.74%  '

$ws.Range("D10").Value = '''0.07979'
$ws.Range("E10").Value = '  +0.25%  '

$ws.Range("D11").Value = '''0.9893'
$ws.Range("E11").Value = '  -0.49%  '

$ws.Range("D12").Value = '''22.43'
$ws.Range("E12").Value = '  +3.64%  '

$ws.Range("D13").Value = '1.897.90'
$ws.Range("E13").Value = '  +0.40%  '

$ws.Range("D14").Value = '''5.838'
$ws.Range("E14").Value = '  -1.18%  '

$ws.Range("D15").Value = '''7.023'
$ws.Range("E15").Value = '  -0.54%  '

$ws.Range("D16").Value = '''1.005'
$ws.Range("E16").Value = '  +0.38%  '

$ws.Range("D17").Value = '''88.85'
$ws.Range("E17").Value = '  +0.46%  '

$ws.Range("D18").Value = '''0.06628'
$ws.Range("E18").Value = '  +0.95%  '

$ws.Range("D19").Value = '''0.00001023'
$ws.Range("E19").Value = '  -0.20%  '

$ws.Range("D20").Value = '''17.44'
$ws.Range("E20").Value = '  +0.30%  '

$ws.Range("D21").Value = '''0.9979'
$ws.Range("E21").Value = '  -0.30%  '

$ws.Range("D22").Value = '29.147.80'
$ws.Range("E22").Value = '  +0.67%  '

$ws.Range("D23").Value = '''5.496'
$ws.Range("E23").Value = '  +1.57%  '

$ws.Range("D24").Value = '''11.42'
$ws.Range("E24").Value = '  -0.16%  '

$ws.Range("D25").Value = '''2.205'
$ws.Range("E25").Value = '  +0.31%  '

$ws.Range("D26").Value = '2.121.10'
$ws.Range("E26").Value = '  +0.28%  '

$ws.Range("D27").Value = '''154.11'
$ws.Range("E27").Value = '  -1.18%  '

$ws.Range("D28").Value = '''19.61'
$ws.Range("E28").Value = '  +0.27%  '

$ws.Range("D29").Value = '''5.990'
$ws.Range("E29").Value = '  +9.43%  '

$ws.Range("D30").Value = '''2.078'
$ws.Range("E30").Value = '  -0.14%  '

$ws.Range("D31").Value = '''116.85'
$ws.Range("E31").Value = '  -0.49%  '

$ws.Range("D32").Value = '''1.053'
$ws.Range("E32").Value = '  +2.05%  '

$ws.Range("D33").Value = '''0.09446'
$ws.Range("E33").Value = '  +1.37%  '

$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = '''1.396'
$ws.Range("E34").Value = '  -0.35%  '

$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '''3.555'
$ws.Range("E35").Value = '  +0.78%  '

$ws.Range("D36").Value = '''5.327'
$ws.Range("E36").Value = '  +0.73%  '

$ws.Range("D37").Value = '''0.06061'
$ws.Range("E37").Value = '  +0.01%  '

$ws.Range("D38").Value = '''0.02238'
$ws.Range("E38").Value = '  +0.53%  '

$ws.Range("D39").Value = '''1.168'
$ws.Range("E39").Value = '  -0.19%  '

$ws.Range("D40").Value = '''8.049'
$ws.Range("E40").Value = '  -3.53%  '

$ws.Range("D41").Value = '''0.5794'
$ws.Range("E41").Value = '  +0.23%  '

$ws.Range("D42").Value = '''0.1820'
$ws.Range("E42").Value = '  -0.09%  '

$ws.Range("D43").Value = '''2.464'
$ws.Range("E43").Value = '  +9.17%  '

$ws.Range("D44").Value = '''10.02'
$ws.Range("E44").Value = '  -0.45%  '

$ws.Range("D45").Value = '''0.07674'
$ws.Range("E45").Value = '  +2.34%  '

$ws.Range("D46").Value = '''1.253'
$ws.Range("E46").Value = '  -0.47%  '

$ws.Range("D47").Value = '''12.01'
$ws.Range("E47").Value = '  +0.75%  '

$ws.Range("D48").Value = '''0.5454'
$ws.Range("E48").Value = '  +0.07%  '

$ws.Range("D49").Value = '''1.894'
$ws.Range("E49").Value = '  -0.25%  '

$ws.Range("D50").Value = '''113.20'
$ws.Range("E50").Value = '  +1.77%  '

$ws.Range("D51").Value = '''43.60'
$ws.Range("E51").Value = '  -2.66%  '
